$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25, pushing existing rows 25-54 down to 26-55.
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new weekly price record.
$ws.Cells.Item(25, 1).Value = 7
$ws.Cells.Item(25, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(25, 3).Value = "Ñuble"
$ws.Cells.Item(25, 4).Value = 44580
$ws.Cells.Item(25, 5).Value = 16
$ws.Cells.Item(25, 6).Value = 100112022
$ws.Cells.Item(25, 7).Value = "Arveja Verde"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 100
$ws.Cells.Item(25, 11).Value = 24000
$ws.Cells.Item(25, 12).Value = 25000
$ws.Cells.Item(25, 13).Value = 24500
$ws.Cells.Item(25, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(25, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(25, 16).Value = 980
$ws.Cells.Item(25, 17).Value = 25
$ws.Cells.Item(25, 18).Value = "Hortaliza"
